# Reorder the comma-separated "Recorded By" names in column G:
# move the first name to the end of the list (left-rotate by one)
# for every row whose G cell contains more than one name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($val -and $val.Contains(",")) {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
            $cell.Value = $rotated -join ", "
        }
    }
}
